# Weekly CompStat update: new crime data collected.
# Updates the report "Volume/Number" and "Week Covering" header strings,
# plus the weekly/28-day/YTD/2-yr/15-yr/32-yr crime-count and percentage
# figures for rows 16-31 (Robbery ... Hate Crimes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header shared strings -------------------------------------------------
$ws.Range("A8").Value  = "Volume 32   Number  31"
$ws.Range("C9").Value  = "Report Covering the Week  7/28/2025  Through  8/3/2025"

# ---- Helper: set a numeric value, fixing up style/number-format if the
# ---- cell previously held a text placeholder ("0" / "***.*") instead of a
# ---- real number. ----------------------------------------------------------
function Set-Num($addr, $val, $fmt) {
    $ws.Range($addr).Value = $val
    if ($fmt) {
        $ws.Range($addr).NumberFormat = $fmt
    }
}

$FMT_INT = "#,##0"
$FMT_PCT = "#,##0.0;`"-`"#,##0.0"

# ---- Row 16: Robbery --------------------------------------------------------
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 71.428571428571
$ws.Range("I16").Value = 63
$ws.Range("J16").Value = 63
$ws.Range("M16").Value = 10.526315789473
$ws.Range("N16").Value = -83.421052631578

# ---- Row 17: Fel. Assault ----------------------------------------------------
$ws.Range("C17").Value = 5
$ws.Range("E17").Value = 66.666666666666
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 50
$ws.Range("I17").Value = 93
$ws.Range("J17").Value = 71
$ws.Range("K17").Value = 30.985915492957
$ws.Range("L17").Value = -10.576923076923
$ws.Range("M17").Value = 9.411764705882
$ws.Range("N17").Value = -37.162162162162

# ---- Row 18: Burglary --------------------------------------------------------
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -33.333333333333
$ws.Range("I18").Value = 62
$ws.Range("J18").Value = 65
$ws.Range("K18").Value = -4.615384615384
$ws.Range("L18").Value = -13.888888888888
$ws.Range("M18").Value = 6.896551724137
$ws.Range("N18").Value = -81.212121212121

# ---- Row 19: Gr. Larceny -----------------------------------------------------
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = -11.764705882352
$ws.Range("F19").Value = 42
$ws.Range("G19").Value = 63
$ws.Range("H19").Value = -33.333333333333
$ws.Range("I19").Value = 356
$ws.Range("J19").Value = 405
$ws.Range("K19").Value = -12.098765432098
$ws.Range("L19").Value = -23.110151187905
$ws.Range("M19").Value = 0.281690140845
$ws.Range("N19").Value = -23.931623931623

# ---- Row 20: G.L.A. ----------------------------------------------------------
# C20 switches from a text placeholder ("0") to a real number.
Set-Num "C20" 4 $FMT_INT
$ws.Range("F20").Value = 5
$ws.Range("H20").Value = 400
$ws.Range("I20").Value = 15
$ws.Range("K20").Value = -44.444444444444
$ws.Range("L20").Value = -59.459459459459
$ws.Range("M20").Value = -40
$ws.Range("N20").Value = -94.773519163763

# ---- Row 21: TOTAL ------------------------------------------------------------
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = 7.407407407407
$ws.Range("F21").Value = 81
$ws.Range("G21").Value = 95
$ws.Range("H21").Value = -14.736842105263
$ws.Range("I21").Value = 601
$ws.Range("J21").Value = 639
$ws.Range("K21").Value = -5.946791862284
$ws.Range("L21").Value = -21.231979030144
$ws.Range("M21").Value = 2.385008517887
$ws.Range("N21").Value = -63.174019607843

# ---- Row 22: Transit -----------------------------------------------------------
# C22/D22/E22 switch from text placeholders ("0"/"0"/"***.*") to real numbers.
Set-Num "C22" 1 $FMT_INT
Set-Num "D22" 1 $FMT_INT
Set-Num "E22" 0 $FMT_PCT
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 150
$ws.Range("I22").Value = 19
$ws.Range("J22").Value = 20
$ws.Range("K22").Value = -5
$ws.Range("L22").Value = -9.523809523809
$ws.Range("M22").Value = 216.666666666667

# ---- Row 23: Housing -----------------------------------------------------------
$ws.Range("F23").Value = 1
$ws.Range("H23").Value = -75
$ws.Range("J23").Value = 27
$ws.Range("K23").Value = -7.407407407407
$ws.Range("L23").Value = -16.666666666666
$ws.Range("M23").Value = -13.793103448275

# ---- Row 24: Petit Larceny -------------------------------------------------------
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = -21.428571428571
$ws.Range("G24").Value = 52
$ws.Range("H24").Value = -3.846153846153
$ws.Range("I24").Value = 489
$ws.Range("J24").Value = 391
$ws.Range("K24").Value = 25.063938618925
$ws.Range("L24").Value = 4.710920770877
$ws.Range("M24").Value = -7.034220532319

# ---- Row 25: Retail Theft --------------------------------------------------------
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 14.285714285714
$ws.Range("F25").Value = 22
$ws.Range("G25").Value = 18
$ws.Range("H25").Value = 22.222222222222
$ws.Range("I25").Value = 267
$ws.Range("J25").Value = 163
$ws.Range("K25").Value = 63.803680981595
$ws.Range("L25").Value = -3.610108303249

# ---- Row 26: Misd. Assault ---------------------------------------------------------
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 233.333333333333
$ws.Range("F26").Value = 31
$ws.Range("G26").Value = 18
$ws.Range("H26").Value = 72.222222222222
$ws.Range("I26").Value = 183
$ws.Range("J26").Value = 185
$ws.Range("K26").Value = -1.081081081081
$ws.Range("L26").Value = -9.405940594059
$ws.Range("M26").Value = -9.852216748768

# ---- Row 28: Other Sex Crimes ---------------------------------------------------------
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = -60
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 34
$ws.Range("J28").Value = 36
$ws.Range("K28").Value = -5.555555555555
$ws.Range("L28").Value = -8.108108108108

# ---- Row 29: Shooting Vic. ---------------------------------------------------------
$ws.Range("M29").Value = -60

# ---- Row 30: Shooting Inc. ---------------------------------------------------------
$ws.Range("M30").Value = -33.333333333333

# ---- Row 31: Hate Crimes ---------------------------------------------------------
# F31 switches from a text placeholder ("0") to a real number.
Set-Num "F31" 1 $FMT_INT
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 4
$ws.Range("K31").Value = -20
$ws.Range("L31").Value = -69.230769230769
